# Sun, Jun 28, 2020  7:04:53 PM
#
# 1) Slide 6's table switches to a different built-in table style.
# 2) The deck's theme colour palette is swapped from the custom "Integral"
#    palette to the stock "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{59D055F3-9BA7-4D96-9250-2B4513B81393}")
    }
}

# --- 2. Swap the theme colour scheme onto the stock "Office Theme" --------
function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Colors(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47",
    "0563C1", "954F72"
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $themeColorScheme.Colors($i + 1).RGB = HexToComRgb $officeThemeColors[$i]
}
